$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric remain stored as text
# (mirrors the inlineStr text cells in the source data).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.261.01"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.648.21"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "217.29"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "0.506"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "0.258"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "0.0637"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "20.06"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.877.31"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.30"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "1.635.23"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "0.554"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "0.0₃0767"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "63.74"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "26.266.56"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "0.999"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "196.44"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "4.45"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "10.05"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "6.35"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").Value = "143.46"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "0.126"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "0.0502"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "3.26"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "0.916"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").Value = "1.142.28"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "0.555"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "5.67"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").Value = "100.53"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").Value = "1.786.30"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "56.27"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0518"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  +5.61%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.418"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0979"
$ws.Range("E51").Value = "  +3.03%  "
